# stm32 cube 1.8.0, LL i2c screen work, leadscrew 2mm
#
# - "feed" sheet: leadscrew pitch (C6) changes from 1.5mm to 2mm.
# - "feed" sheet: F9/G9 (micro-stepping hint columns for the first,
#   smallest step entry) are cleared out.
# - "feed" sheet: selection moves from I17 to A20.
# - "ramp up" sheet: rows 9:35 get their explicit row height cleared
#   (autofit back to the sheet default).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("feed")

# Leadscrew pitch: 1.5 -> 2 mm. All the dependent formula cells
# (K3, J4, K4, J6, B9:B40, I9:I40, ...) recalc automatically.
$ws.Range("C6").Value = 2

# Clear the stray hint values that used to live in F9/G9.
$ws.Range("F9:G9").ClearContents()

# Selection moved from I17 to A20.
$ws.Range("A20").Select()

$ws3 = $wb.Worksheets.Item("ramp up")

# Autofit rows 9:35 back to the sheet's default row height.
$ws3.Rows("9:35").AutoFit()
